$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the bit values on row 2 (register 0x01 / ID)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0

# Update the bit values on row 14 (register 0x0D / LOFF_SENSP)
$ws.Range("J14").Value = 1

# Update the bit values on row 15 (register 0x0E / LOFF_SENSN)
$ws.Range("J15").Value = 1

# Recalculate formulas so the derived hex text values (C2, C14, C15) refresh
$excel.Calculate()

# Update the active cell selection to match the saved view state
$ws.Range("F2").Select()
